$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: A22 changes from a blank " " text marker to the plain number 21
$ws.Range("A22").Value = 21

# Rows 23-28: new data rows -- running index in column A, and a little
# A/B/C/D/E/F legend filled into columns K:P for each row.
$rows = 23..28
$seq  = 22
foreach ($r in $rows) {
    $ws.Range("A$r").Value = $seq
    $ws.Range("K$r").Value = "A"
    $ws.Range("L$r").Value = "B"
    $ws.Range("M$r").Value = "C"
    $ws.Range("N$r").Value = "D"
    $ws.Range("O$r").Value = "E"
    $ws.Range("P$r").Value = "F"
    $seq++
}

# Restore the saved selection/active cell state
$ws.Range("A20:A28").Select()
